$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Swerve module swap out recalibrated the hood-angle LUT entry for this
# distance sample: update the measured distance in A17. The dependent
# formulas in B17 (predicted hood angle) and C17 (predicted velocity)
# recalc automatically.
$ws.Range("A17").Value = 11.913
